$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 5: "Top K Frequent Elements" (number frequency with hash) entry
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 347
$ws.Range("C5").Value = "给定一个非空的整数数组，返回其中出现频率前 k 高的元素。 "
$ws.Range("D5").Value = "1 hashmap存储每个数字对应的出现次数`n2 创建小顶堆，小顶堆的长度是k，存放的顺序是比较这个数字出现的次数。`n3 先添加元素【add】，如果堆长度大于k，就移除堆头部元素即出现次数最小的元素【remove/poll】；长度小于k，就继续循环。`n4 所有数字添加完成，小顶堆也就创建完成"
$ws.Range("E5").Value = "小顶堆`n哈希表"
$ws.Range("G5").Value = "O(n)`n堆的长度是k，哈希表的长度是n，"
$ws.Range("F5").Value = "O(nlogk)`nn是数组长度`nk是堆的长度"

# Match the row height Excel computed after wrapping the new multi-line content
$ws.Rows.Item(5).RowHeight = 154

# Selection moved to E7 in the saved view
$ws.Range("E7").Select()
